$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 5667.4062
$ws.Range("I15").Value = 5667.4062
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 17002.2186
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -16833.2186

$ws.Range("H28").Value = 559.5
$ws.Range("I28").Value = 459
$ws.Range("J28").Value = 660
$ws.Range("K28").Value = 459
$ws.Range("L28").Value = 660
$ws.Range("M28").Value = 26
$ws.Range("N28").Value = -1630

$ws.Range("H129").Value = 1105.6129
$ws.Range("I129").Value = 2095.6667
$ws.Range("J129").Value = 868
$ws.Range("K129").Value = 6287.000100000001
$ws.Range("L129").Value = 2604
$ws.Range("M129").Value = -1287.000100000001
$ws.Range("N129").Value = -12604

$ws.Range("H132").Value = 39156.27
$ws.Range("I132").Value = 25946.72
$ws.Range("J132").Value = 102268.555
$ws.Range("K132").Value = 77840.16
$ws.Range("L132").Value = 306805.665
$ws.Range("M132").Value = -75310.16
$ws.Range("N132").Value = -311865.665

$ws.Range("H137").Value = 2144133
$ws.Range("I137").Value = 4281166
$ws.Range("J137").Value = 7100.1113
$ws.Range("K137").Value = 12843498
$ws.Range("L137").Value = 21300.3339
$ws.Range("M137").Value = -12840948
$ws.Range("N137").Value = -26400.3339

$ws.Range("H138").Value = 2599.6667
$ws.Range("I138").Value = 2163.7778
$ws.Range("J138").Value = 2817.611
$ws.Range("K138").Value = 6491.3334
$ws.Range("L138").Value = 8452.832999999999
$ws.Range("M138").Value = -1351.3334
$ws.Range("N138").Value = -18732.833

$ws.Range("H141").Value = 2277.6667
$ws.Range("I141").Value = 1349.5
$ws.Range("J141").Value = 5526.25
$ws.Range("K141").Value = 4048.5
$ws.Range("L141").Value = 16578.75
$ws.Range("M141").Value = 1131.5
$ws.Range("N141").Value = -26938.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3071
$ws.Range("I61").Value = 2584.6
$ws.Range("J61").Value = 3395.2666
$ws.Range("K61").Value = 2584.6
$ws.Range("L61").Value = 3395.2666
$ws.Range("M61").Value = -2372.6
$ws.Range("N61").Value = -3819.2666

$ws.Range("H74").Value = 1862.5667
$ws.Range("I74").Value = 1522.1818
$ws.Range("J74").Value = 2798.625
$ws.Range("K74").Value = 1522.1818
$ws.Range("L74").Value = 2798.625
$ws.Range("M74").Value = -648.1818000000001
$ws.Range("N74").Value = -4546.625

$ws.Range("H77").Value = 1862.5667
$ws.Range("I77").Value = 1522.1818
$ws.Range("J77").Value = 2798.625
$ws.Range("K77").Value = 7610.909000000001
$ws.Range("L77").Value = 13993.125
$ws.Range("M77").Value = -3242.909000000001
$ws.Range("N77").Value = -22729.125

$ws.Range("H132").Value = 2096.9531
$ws.Range("I132").Value = 1564.3726
$ws.Range("J132").Value = 4186.3076
$ws.Range("K132").Value = 4693.1178
$ws.Range("L132").Value = 12558.9228
$ws.Range("M132").Value = -2163.1178
$ws.Range("N132").Value = -17618.9228

$ws.Range("H136").Value = 3071
$ws.Range("I136").Value = 2584.6
$ws.Range("J136").Value = 3395.2666
$ws.Range("K136").Value = 7753.799999999999
$ws.Range("L136").Value = 10185.7998
$ws.Range("M136").Value = -5203.799999999999
$ws.Range("N136").Value = -15285.7998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2645.261
$ws.Range("I105").Value = 2170
$ws.Range("J105").Value = 3263.1
$ws.Range("K105").Value = 2170
$ws.Range("L105").Value = 3263.1
$ws.Range("M105").Value = -423
$ws.Range("N105").Value = -6757.1

$ws.Range("H107").Value = 1999.7241
$ws.Range("I107").Value = 1441.4
$ws.Range("J107").Value = 3240.4443
$ws.Range("K107").Value = 1441.4
$ws.Range("L107").Value = 3240.4443
$ws.Range("M107").Value = 478.5999999999999
$ws.Range("N107").Value = -7080.4443

$ws.Range("H134").Value = 1919.3077
$ws.Range("I134").Value = 1566.238
$ws.Range("J134").Value = 3402.2
$ws.Range("K134").Value = 4698.714
$ws.Range("L134").Value = 10206.6
$ws.Range("M134").Value = -2163.714
$ws.Range("N134").Value = -15276.6

$ws.Range("H135").Value = 55242.777
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 55242.777
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 55242.777
$ws.Range("N135").Value = -65382.777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1789.6428
$ws.Range("I22").Value = 331.875
$ws.Range("J22").Value = 3733.3333
$ws.Range("K22").Value = 331.875
$ws.Range("L22").Value = 3733.3333
$ws.Range("M22").Value = 18.125
$ws.Range("N22").Value = -4433.3333

$ws.Range("H31").Value = 4277088
$ws.Range("I31").Value = 1446.7941
$ws.Range("J31").Value = 7580993
$ws.Range("K31").Value = 1446.7941
$ws.Range("L31").Value = 7580993
$ws.Range("M31").Value = -1151.7941
$ws.Range("N31").Value = -7581583

$ws.Range("H34").Value = 4277088
$ws.Range("I34").Value = 1446.7941
$ws.Range("J34").Value = 7580993
$ws.Range("K34").Value = 1446.7941
$ws.Range("L34").Value = 7580993
$ws.Range("M34").Value = -1244.7941
$ws.Range("N34").Value = -7581397

$ws.Range("H51").Value = 71461380
$ws.Range("I51").Value = 500000000
$ws.Range("J51").Value = 38274.832
$ws.Range("K51").Value = 500000000
$ws.Range("L51").Value = 38274.832
$ws.Range("M51").Value = -499999264
$ws.Range("N51").Value = -39746.832

$ws.Range("H58").Value = 1804.8387
$ws.Range("I58").Value = 1040.85
$ws.Range("J58").Value = 3193.9092
$ws.Range("K58").Value = 1040.85
$ws.Range("L58").Value = 3193.9092
$ws.Range("M58").Value = -837.8499999999999
$ws.Range("N58").Value = -3599.9092

$ws.Range("H61").Value = 71461380
$ws.Range("I61").Value = 500000000
$ws.Range("J61").Value = 38274.832
$ws.Range("K61").Value = 500000000
$ws.Range("L61").Value = 38274.832
$ws.Range("M61").Value = -499999652
$ws.Range("N61").Value = -38970.832

$ws.Range("H132").Value = 47816.13
$ws.Range("I132").Value = 2132.318
$ws.Range("J132").Value = 159487.67
$ws.Range("K132").Value = 6396.954000000001
$ws.Range("L132").Value = 478463.01
$ws.Range("M132").Value = -3866.954000000001
$ws.Range("N132").Value = -483523.01

$ws.Range("H133").Value = 32450
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 32450
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 32450
$ws.Range("N133").Value = -37510

$ws.Range("H134").Value = 1081558.6
$ws.Range("I134").Value = 613501.2
$ws.Range("J134").Value = 4669999
$ws.Range("K134").Value = 1840503.6
$ws.Range("L134").Value = 14009997
$ws.Range("M134").Value = -1837968.6
$ws.Range("N134").Value = -14015067

$ws.Range("H136").Value = 1804.8387
$ws.Range("I136").Value = 1040.85
$ws.Range("J136").Value = 3193.9092
$ws.Range("K136").Value = 3122.55
$ws.Range("L136").Value = 9581.7276
$ws.Range("M136").Value = -572.5499999999997
$ws.Range("N136").Value = -14681.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 100333.336
$ws.Range("I130").Value = 100333.336
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 301000.008
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -295980.008
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value = 2311.5334
$ws.Range("I132").Value = 973.25
$ws.Range("J132").Value = 2798.182
$ws.Range("K132").Value = 8759.25
$ws.Range("L132").Value = 25183.638
$ws.Range("M132").Value = -6229.25
$ws.Range("N132").Value = -30243.638

$ws.Range("H134").Value = 47764424
$ws.Range("I134").Value = 50152292
$ws.Range("J134").Value = 7070
$ws.Range("K134").Value = 150456876
$ws.Range("L134").Value = 21210
$ws.Range("M134").Value = -150451806
$ws.Range("N134").Value = -31350

$ws.Range("H139").Value = 184355.3
$ws.Range("I139").Value = 184355.3
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 553065.8999999999
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -547925.8999999999
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H132").Value = 22729584
$ws.Range("I132").Value = 35715972
$ws.Range("J132").Value = 3402.0625
$ws.Range("K132").Value = 107147916
$ws.Range("L132").Value = 10206.1875
$ws.Range("M132").Value = -107145386
$ws.Range("N132").Value = -15266.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 4718.143
$ws.Range("I132").Value = 3544.818
$ws.Range("J132").Value = 6008.8
$ws.Range("K132").Value = 10634.454
$ws.Range("L132").Value = 18026.4
$ws.Range("M132").Value = -8104.454000000002
$ws.Range("N132").Value = -23086.4

$ws.Range("H136").Value = 1722.3
$ws.Range("I136").Value = 1307.1154
$ws.Range("J136").Value = 4421
$ws.Range("K136").Value = 3921.3462
$ws.Range("L136").Value = 13263
$ws.Range("M136").Value = -1371.3462
$ws.Range("N136").Value = -18363

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1451239.9
$ws.Range("I132").Value = 1977531.8
$ws.Range("J132").Value = 3937
$ws.Range("K132").Value = 5932595.4
$ws.Range("L132").Value = 11811
$ws.Range("M132").Value = -5930065.4
$ws.Range("N132").Value = -16871

$ws.Range("H136").Value = 778791.25
$ws.Range("I136").Value = 898263
$ws.Range("J136").Value = 2224.75
$ws.Range("K136").Value = 2694789
$ws.Range("L136").Value = 6674.25
$ws.Range("M136").Value = -2692239
$ws.Range("N136").Value = -11774.25
